$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.5260084838407408
$ws.Range("C2").Value = 0.2872063558142699
$ws.Range("D2").Value = 0.05618599683036507
$ws.Range("F2").Value = 1.123251065952957
$ws.Range("G2").Value = 0.002480392340157849
$ws.Range("K2").Value = 0.2451419974972282
$ws.Range("L2").Value = 0.2941846220858366
$ws.Range("M2").Value = 0.1852605062729822
$ws.Range("O2").Value = 4.067591855727045

$ws.Range("B3").Value = 0.4909891331739971
$ws.Range("C3").Value = 0.2881809373423287
$ws.Range("D3").Value = 0.05441949097095033
$ws.Range("F3").Value = 1.123706267813873
$ws.Range("G3").Value = 0.002482947316076259
$ws.Range("K3").Value = 0.2148992048153957
$ws.Range("L3").Value = 0.2908679756148871
$ws.Range("M3").Value = 0.1781141232952024
$ws.Range("O3").Value = 4.083599013994387

$ws.Range("B4").Value = 0.4696698222861926
$ws.Range("C4").Value = 0.2888166716796761
$ws.Range("D4").Value = 0.05332109844982824
$ws.Range("F4").Value = 1.124571423672705
$ws.Range("G4").Value = 0.002484600572380495
$ws.Range("K4").Value = 0.1963093341437769
$ws.Range("L4").Value = 0.2889787838425733
$ws.Range("M4").Value = 0.1738111013135963
$ws.Range("O4").Value = 4.09538613175738

$ws.Range("B5").Value = 0.4610285048303808
$ws.Range("C5").Value = 0.2890851550677738
$ws.Range("D5").Value = 0.05287005678045631
$ws.Range("F5").Value = 1.125071314183174
$ws.Range("G5").Value = 0.002485295596443295
$ws.Range("K5").Value = 0.1887289998868624
$ws.Range("L5").Value = 0.2882460103439612
$ws.Range("M5").Value = 0.1720790392824263
$ws.Range("O5").Value = 4.100682132642248

$ws.Range("B6").Value = 0.4595964432696178
$ws.Range("C6").Value = 0.2891303060801285
$ws.Range("D6").Value = 0.05279495460925432
$ws.Range("F6").Value = 1.125163221011064
$ws.Range("G6").Value = 0.002485412293665361
$ws.Range("K6").Value = 0.187470012033117
$ws.Range("L6").Value = 0.2881265756887572
$ws.Range("M6").Value = 0.1717927307103224
$ws.Range("O6").Value = 4.101591288659307

$ws.Range("B7").Value = 0.469553093577872
$ws.Range("C7").Value = 0.2888202543762937
$ws.Range("D7").Value = 0.05331502943749911
$ws.Range("F7").Value = 1.124577568744712
$ws.Range("G7").Value = 0.002484609859490373
$ws.Range("K7").Value = 0.1962071220695378
$ws.Range("L7").Value = 0.2889687511542505
$ws.Range("M7").Value = 0.1737876551063593
$ws.Range("O7").Value = 4.095455560565455

$ws.Range("B8").Value = 0.5138961959052892
$ws.Range("C8").Value = 0.2875346587847964
$ws.Range("D8").Value = 0.05557977105444678
$ws.Range("F8").Value = 1.123286501178328
$ws.Range("G8").Value = 0.002481255797635958
$ws.Range("K8").Value = 0.2347187962640618
$ws.Range("L8").Value = 0.293010517712645
$ws.Range("M8").Value = 0.1827788780771655
$ws.Range("O8").Value = 4.072704679758488

$ws.Range("B9").Value = 0.6022840684067035
$ws.Range("C9").Value = 0.2853086249970254
$ws.Range("D9").Value = 0.05991112569153501
$ws.Range("F9").Value = 1.12540021635651
$ws.Range("G9").Value = 0.002475346020455953
$ws.Range("K9").Value = 0.3100627074197462
$ws.Range("L9").Value = 0.3021027377546233
$ws.Range("M9").Value = 0.2010804578281267
$ws.Range("O9").Value = 4.043628215787692

$ws.Range("B10").Value = 0.6680771914600996
$ws.Range("C10").Value = 0.2838513065478061
$ws.Range("D10").Value = 0.06302578251229107
$ws.Range("F10").Value = 1.129785203476004
$ws.Range("G10").Value = 0.002471407029238181
$ws.Range("K10").Value = 0.3652972908333822
$ws.Range("L10").Value = 0.3094922771653756
$ws.Range("M10").Value = 0.214931567319482
$ws.Range("O10").Value = 4.031737999728051

$ws.Range("B11").Value = 0.6981903876179558
$ws.Range("C11").Value = 0.283226660758011
$ws.Range("D11").Value = 0.06442792264946462
$ws.Range("F11").Value = 1.132395131648209
$ws.Range("G11").Value = 0.002469701716072035
$ws.Range("K11").Value = 0.390396455444062
$ws.Range("L11").Value = 0.3130077718672624
$ws.Range("M11").Value = 0.2213201031223235
$ws.Range("O11").Value = 4.028385937094896

$ws.Range("B12").Value = 0.7096194320766642
$ws.Range("C12").Value = 0.2829956032134895
$ws.Range("D12").Value = 0.06495674152503028
$ws.Range("F12").Value = 1.133471877417222
$ws.Range("G12").Value = 0.00246906833994529
$ws.Range("K12").Value = 0.3998966116253087
$ws.Range("L12").Value = 0.314361091532902
$ws.Range("M12").Value = 0.2237517843034169
$ws.Range("O12").Value = 4.027412331638516

$ws.Range("B13").Value = 0.7071568430538662
$ws.Range("C13").Value = 0.2830451221533465
$ws.Range("D13").Value = 0.06484294650548605
$ws.Range("F13").Value = 1.133236048960697
$ws.Range("G13").Value = 0.002469204198803517
$ws.Range("K13").Value = 0.3978507833287779
$ws.Range("L13").Value = 0.3140686489657583
$ws.Range("M13").Value = 0.223227524972387
$ws.Range("O13").Value = 4.027608861944799

$ws.Range("B14").Value = 0.6991301475441958
$ws.Range("C14").Value = 0.2832075417821898
$ws.Range("D14").Value = 0.06447147193759406
$ws.Range("F14").Value = 1.132481944276208
$ws.Range("G14").Value = 0.002469649359913812
$ws.Range("K14").Value = 0.3911781286663256
$ws.Range("L14").Value = 0.3131186681600582
$ws.Range("M14").Value = 0.2215199095105191
$ws.Range("O14").Value = 4.028299910910192

$ws.Range("B15").Value = 0.6942169118765094
$ws.Range("C15").Value = 0.2833077417207264
$ws.Range("D15").Value = 0.06424365352624761
$ws.Range("F15").Value = 1.13203154744545
$ws.Range("G15").Value = 0.002469923645583137
$ws.Range("K15").Value = 0.3870903556183407
$ws.Range("L15").Value = 0.3125396509077092
$ws.Range("M15").Value = 0.2204755677596566
$ws.Range("O15").Value = 4.028761712683746

$ws.Range("B16").Value = 0.6661128687456426
$ws.Range("C16").Value = 0.2838928971135175
$ws.Range("D16").Value = 0.0629338509604338
$ws.Range("F16").Value = 1.129627014136148
$ws.Range("G16").Value = 0.002471520212428544
$ws.Range("K16").Value = 0.3636564153241579
$ws.Range("L16").Value = 0.3092656241351079
$ws.Range("M16").Value = 0.2145158143188937
$ws.Range("O16").Value = 4.031998451959652

$ws.Range("B17").Value = 0.6489185624289746
$ws.Range("C17").Value = 0.2842616618191478
$ws.Range("D17").Value = 0.06212653996337281
$ws.Range("F17").Value = 1.128309444708393
$ws.Range("G17").Value = 0.002472521782107438
$ws.Range("K17").Value = 0.349273151728454
$ws.Range("L17").Value = 0.3072965035666613
$ws.Range("M17").Value = 0.2108820523888895
$ws.Range("O17").Value = 4.03451087916585

$ws.Range("B18").Value = 0.6390461649270094
$ws.Range("C18").Value = 0.2844773715443587
$ws.Range("D18").Value = 0.06166081064765905
$ws.Range("F18").Value = 1.12760951923913
$ws.Range("G18").Value = 0.002473106008784321
$ws.Range("K18").Value = 0.3409977322494626
$ws.Range("L18").Value = 0.3061784122143507
$ws.Range("M18").Value = 0.2088002605129589
$ws.Range("O18").Value = 4.036149563476272

$ws.Range("B19").Value = 0.6357065315440309
$ws.Range("C19").Value = 0.2845510272437117
$ws.Range("D19").Value = 0.0615028854533719
$ws.Range("F19").Value = 1.127382483017641
$ws.Range("G19").Value = 0.00247330521926433
$ws.Range("K19").Value = 0.3381953929388146
$ws.Range("L19").Value = 0.3058023370184344
$ws.Range("M19").Value = 0.2080968223258424
$ws.Range("O19").Value = 4.036737646033174

$ws.Range("B20").Value = 0.6507471380009235
$ws.Range("C20").Value = 0.2842220331553627
$ws.Range("D20").Value = 0.06221262314308262
$ws.Range("F20").Value = 1.128443709486191
$ws.Range("G20").Value = 0.002472414320015437
$ws.Range("K20").Value = 0.350804541337709
$ws.Range("L20").Value = 0.3075046201025771
$ws.Range("M20").Value = 0.2112680193448071
$ws.Range("O20").Value = 4.034223390028899

$ws.Range("B21").Value = 0.7014870860635369
$ws.Range("C21").Value = 0.2831596866074939
$ws.Range("D21").Value = 0.06458064131297903
$ws.Range("F21").Value = 1.132701043652403
$ws.Range("G21").Value = 0.002469518269483572
$ws.Range("K21").Value = 0.3931381705883723
$ws.Range("L21").Value = 0.3133971018405077
$ws.Range("M21").Value = 0.2220211397412299
$ws.Range("O21").Value = 4.028088906773974

$ws.Range("B22").Value = 0.7347988736070477
$ws.Range("C22").Value = 0.2824973253387562
$ws.Range("D22").Value = 0.06611578731004641
$ws.Range("F22").Value = 1.135998833531758
$ws.Range("G22").Value = 0.002467697719614843
$ws.Range("K22").Value = 0.4207800779605293
$ws.Range("L22").Value = 0.3173768392162231
$ws.Range("M22").Value = 0.2291216229365176
$ws.Range("O22").Value = 4.025803499780949

$ws.Range("B23").Value = 0.7170061804222883
$ws.Range("C23").Value = 0.2828479252179719
$ws.Range("D23").Value = 0.06529760134992557
$ws.Range("F23").Value = 1.134191597192299
$ws.Range("G23").Value = 0.002468662795073282
$ws.Range("K23").Value = 0.4060295506720308
$ws.Range("L23").Value = 0.3152410261526057
$ws.Range("M23").Value = 0.2253253475720598
$ws.Range("O23").Value = 4.026865542491294

$ws.Range("B24").Value = 0.6499203989538387
$ws.Range("C24").Value = 0.2842399377406188
$ws.Range("D24").Value = 0.06217370990568583
$ws.Range("F24").Value = 1.128382829071839
$ws.Range("G24").Value = 0.002472462877630591
$ws.Range("K24").Value = 0.3501122196561255
$ws.Range("L24").Value = 0.3074104870599257
$ws.Range("M24").Value = 0.2110935009177766
$ws.Range("O24").Value = 4.034352758761514

$ws.Range("B25").Value = 0.5782214489091757
$ws.Range("C25").Value = 0.2858794197659691
$ws.Range("D25").Value = 0.0587512020104981
$ws.Range("F25").Value = 1.124331062842501
$ws.Range("G25").Value = 0.002476873725734377
$ws.Range("K25").Value = 0.2897004658025537
$ws.Range("L25").Value = 0.2995183465979778
$ws.Range("M25").Value = 0.1960580389806879
$ws.Range("O25").Value = 4.049830732091493

